# Intégrer 195 propositions pour 11 candidats dans 8 villes
# Updates "Programmes candidats" (candidate rows) and "Resume par ville"
# (per-city summary rows) to reflect newly-found full/partial programmes.

$wb = $excel.ActiveWorkbook
$wsCandidats = $wb.Worksheets.Item("Programmes candidats")
$wsVilles    = $wb.Worksheets.Item("Resume par ville")

# ---------------------------------------------------------------------
# Reference cells whose formatting ("style") we reuse via copy/paste of
# formats only, so the existing shared cellXfs entries (colour-coded
# status pills) are reused instead of inventing new ones.
#   Programmes candidats  F24 -> "Programme complet" (green)
#                          F2  -> "Bien couvert"      (yellow)
#   Resume par ville      F5  -> "Programme(s) complet(s)" (green)
#                          F4  -> "Bien couvert"            (yellow)
#                          F2  -> "Partiellement couvert"   (orange)
# ---------------------------------------------------------------------
$progCompletFmt = $wsCandidats.Range("F24")
$bienCouvertFmt = $wsCandidats.Range("F2")

$villeProgCompletFmt  = $wsVilles.Range("F5")
$villeBienCouvertFmt  = $wsVilles.Range("F4")
$villePartielFmt      = $wsVilles.Range("F2")

function Set-StatutCellFormat($range, $formatSource) {
    $formatSource.Copy()
    $range.PasteSpecial(-4122)  # xlPasteFormats
}

# ---------------------------------------------------------------------
# Sheet "Programmes candidats": one row per candidate
# ---------------------------------------------------------------------

# Row 29 - Annecy / Antoine Armand -> Programme complet
$wsCandidats.Range("E29").Value = "https://acteursannecy.fr/programme"
$wsCandidats.Range("F29").Value = "Programme complet"
Set-StatutCellFormat $wsCandidats.Range("F29") $progCompletFmt
$wsCandidats.Range("G29").Value = 25
$wsCandidats.Range("H29").Value = "Oui"

# Row 32 - Annecy / Guillaume Roit-Lévêque -> Programme complet
$wsCandidats.Range("E32").Value = "https://retrouvons-annecy.fr/le-programme/"
$wsCandidats.Range("F32").Value = "Programme complet"
Set-StatutCellFormat $wsCandidats.Range("F32") $progCompletFmt
$wsCandidats.Range("G32").Value = 21
$wsCandidats.Range("H32").Value = "Oui"

# Row 42 - Antony / Perrine Precetti -> Bien couvert
$wsCandidats.Range("E42").Value = "https://www.antonyavenir.fr/programme"
$wsCandidats.Range("F42").Value = "Bien couvert"
Set-StatutCellFormat $wsCandidats.Range("F42") $bienCouvertFmt
$wsCandidats.Range("G42").Value = 10

# Row 123 - Chambéry / Brice Bernard -> Programme complet
$wsCandidats.Range("E123").Value = "https://brice-bernard.fr/programme.html"
$wsCandidats.Range("F123").Value = "Programme complet"
Set-StatutCellFormat $wsCandidats.Range("F123") $progCompletFmt
$wsCandidats.Range("G123").Value = 27
$wsCandidats.Range("H123").Value = "Oui"

# Row 178 - Fréjus / Emmanuel Bonnemain -> Bien couvert
$wsCandidats.Range("F178").Value = "Bien couvert"
Set-StatutCellFormat $wsCandidats.Range("F178") $bienCouvertFmt
$wsCandidats.Range("G178").Value = 14

# Row 193 - Issy-les-Moulineaux / Mathieu Morel -> Programme complet
$wsCandidats.Range("E193").Value = "https://issyecoloetsocial.fr/notre-programme/"
$wsCandidats.Range("F193").Value = "Programme complet"
Set-StatutCellFormat $wsCandidats.Range("F193") $progCompletFmt
$wsCandidats.Range("G193").Value = 21
$wsCandidats.Range("H193").Value = "Oui"

# Row 204 - La Seyne-sur-Mer / Joseph Minniti -> Bien couvert
$wsCandidats.Range("F204").Value = "Bien couvert"
Set-StatutCellFormat $wsCandidats.Range("F204") $bienCouvertFmt
$wsCandidats.Range("G204").Value = 19

# Row 205 - La Seyne-sur-Mer / Cheikh Mansour -> Bien couvert
$wsCandidats.Range("F205").Value = "Bien couvert"
Set-StatutCellFormat $wsCandidats.Range("F205") $bienCouvertFmt
$wsCandidats.Range("G205").Value = 14

# Row 429 - Rueil-Malmaison / Patrick Indjian -> Bien couvert
$wsCandidats.Range("F429").Value = "Bien couvert"
Set-StatutCellFormat $wsCandidats.Range("F429") $bienCouvertFmt
$wsCandidats.Range("G429").Value = 19

# Row 538 - Vitry-sur-Seine / Pierre Bell-Lloch -> Bien couvert
$wsCandidats.Range("E538").Value = "https://pbl2026.fr/programme.html"
$wsCandidats.Range("F538").Value = "Bien couvert"
Set-StatutCellFormat $wsCandidats.Range("F538") $bienCouvertFmt
$wsCandidats.Range("G538").Value = 13

# Row 539 - Vitry-sur-Seine / Hocine Tmimi -> Bien couvert
$wsCandidats.Range("F539").Value = "Bien couvert"
Set-StatutCellFormat $wsCandidats.Range("F539") $bienCouvertFmt
$wsCandidats.Range("G539").Value = 22

# Row 548 - Évry-Courcouronnes / Farida Amrani -> Bien couvert
$wsCandidats.Range("F548").Value = "Bien couvert"
Set-StatutCellFormat $wsCandidats.Range("F548") $bienCouvertFmt
$wsCandidats.Range("G548").Value = 10

$excel.CutCopyMode = $false

# ---------------------------------------------------------------------
# Sheet "Resume par ville": one row per city, aggregated totals
# ---------------------------------------------------------------------

# Row 6 - Annecy -> Programme(s) complet(s)
$wsVilles.Range("C6").Value = 46
$wsVilles.Range("E6").Value = 2
$wsVilles.Range("F6").Value = "Programme(s) complet(s)"
Set-StatutCellFormat $wsVilles.Range("F6") $villeProgCompletFmt

# Row 8 - Antony -> Partiellement couvert
$wsVilles.Range("C8").Value = 10
$wsVilles.Range("F8").Value = "Partiellement couvert"
Set-StatutCellFormat $wsVilles.Range("F8") $villePartielFmt

# Row 25 - Chambéry -> Programme(s) complet(s)
$wsVilles.Range("C25").Value = 27
$wsVilles.Range("E25").Value = 1
$wsVilles.Range("F25").Value = "Programme(s) complet(s)"
Set-StatutCellFormat $wsVilles.Range("F25") $villeProgCompletFmt

# Row 38 - Fréjus -> Partiellement couvert
$wsVilles.Range("C38").Value = 14
$wsVilles.Range("F38").Value = "Partiellement couvert"
Set-StatutCellFormat $wsVilles.Range("F38") $villePartielFmt

# Row 40 - Issy-les-Moulineaux -> Programme(s) complet(s)
$wsVilles.Range("C40").Value = 21
$wsVilles.Range("E40").Value = 1
$wsVilles.Range("F40").Value = "Programme(s) complet(s)"
Set-StatutCellFormat $wsVilles.Range("F40") $villeProgCompletFmt

# Row 43 - La Seyne-sur-Mer -> Bien couvert
$wsVilles.Range("C43").Value = 33
$wsVilles.Range("F43").Value = "Bien couvert"
Set-StatutCellFormat $wsVilles.Range("F43") $villeBienCouvertFmt

# Row 82 - Rueil-Malmaison -> Partiellement couvert
$wsVilles.Range("C82").Value = 19
$wsVilles.Range("F82").Value = "Partiellement couvert"
Set-StatutCellFormat $wsVilles.Range("F82") $villePartielFmt

# Row 101 - Vitry-sur-Seine -> Bien couvert
$wsVilles.Range("C101").Value = 35
$wsVilles.Range("F101").Value = "Bien couvert"
Set-StatutCellFormat $wsVilles.Range("F101") $villeBienCouvertFmt

# Row 104 - Évry-Courcouronnes -> Partiellement couvert
$wsVilles.Range("C104").Value = 10
$wsVilles.Range("F104").Value = "Partiellement couvert"
Set-StatutCellFormat $wsVilles.Range("F104") $villePartielFmt

$excel.CutCopyMode = $false
